$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (for "Coliflor" at "Macroferia Regional de Talca") was recorded.
# It belongs right before the current row 401, so insert a blank row there;
# this pushes the former rows 401-495 down to 402-496.
$ws.Rows(401).Insert()

# Fill in the newly inserted row 401 with its data.
$ws.Cells.Item(401, 1).Value = 5
$ws.Cells.Item(401, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(401, 3).Value = "Maule"
$ws.Cells.Item(401, 4).Value = 45173
$ws.Cells.Item(401, 5).Value = 7
$ws.Cells.Item(401, 6).Value = 100112008
$ws.Cells.Item(401, 7).Value = "Coliflor"
$ws.Cells.Item(401, 8).Value = "Sin especificar"
$ws.Cells.Item(401, 9).Value = "Primera"
$ws.Cells.Item(401, 10).Value = 3000
$ws.Cells.Item(401, 11).Value = 900
$ws.Cells.Item(401, 12).Value = 900
$ws.Cells.Item(401, 13).Value = 900
$ws.Cells.Item(401, 14).Value = '$/unidad'
$ws.Cells.Item(401, 15).Value = "Región del Maule"
$ws.Cells.Item(401, 16).Value = 900
$ws.Cells.Item(401, 17).Value = 1
$ws.Cells.Item(401, 18).Value = "Hortaliza"
